$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 68.392882
$ws.Range("H2").Value = 205.178646
$ws.Range("I2").Value = 0.3817002623156464
$ws.Range("J2").Value = 0.3817002623156463
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.683873666666667
$ws.Range("N2").Value = 8.051621
$ws.Range("O2").Value = 0.7025487437947114
$ws.Range("P2").Value = 0.7025487437947114
$ws.Range("Q2").Value = 183.5578549872407
$ws.Range("R2").Value = 1652.020694885166
$ws.Range("S2").Value = 0.2681630397959692
$ws.Range("T2").Value = 0.2681630397959691

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 68.392882
$ws.Range("H3").Value = 205.178646
$ws.Range("I3").Value = 0.3817002623156464
$ws.Range("J3").Value = 0.3817002623156463
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.5528256666666667
$ws.Range("N3").Value = 1.658477
$ws.Range("O3").Value = 0.1447113485548341
$ws.Range("P3").Value = 0.1447113485548341
$ws.Range("Q3").Value = 37.80934058690467
$ws.Range("R3").Value = 340.284065282142
$ws.Range("S3").Value = 0.05523635970343109
$ws.Range("T3").Value = 0.05523635970343108

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 68.392882
$ws.Range("H4").Value = 205.178646
$ws.Range("I4").Value = 0.3817002623156464
$ws.Range("J4").Value = 0.3817002623156463
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5834963333333333
$ws.Range("N4").Value = 1.750489
$ws.Range("O4").Value = 0.1527399076504545
$ws.Range("P4").Value = 0.1527399076504546
$ws.Range("Q4").Value = 39.90699587309933
$ws.Range("R4").Value = 359.1629628578939
$ws.Range("S4").Value = 0.05830086281624609
$ws.Range("T4").Value = 0.05830086281624609

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 60.20577233333334
$ws.Range("H5").Value = 180.617317
$ws.Range("I5").Value = 0.3360080526004068
$ws.Range("J5").Value = 0.3360080526004068
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.683873666666667
$ws.Range("N5").Value = 8.051621
$ws.Range("O5").Value = 0.7025487437947114
$ws.Range("P5").Value = 0.7025487437947114
$ws.Range("Q5").Value = 161.5846869467619
$ws.Range("R5").Value = 1454.262182520857
$ws.Range("S5").Value = 0.2360620352593231
$ws.Range("T5").Value = 0.2360620352593231

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 60.20577233333334
$ws.Range("H6").Value = 180.617317
$ws.Range("I6").Value = 0.3360080526004068
$ws.Range("J6").Value = 0.3360080526004068
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.5528256666666667
$ws.Range("N6").Value = 1.658477
$ws.Range("O6").Value = 0.1447113485548341
$ws.Range("P6").Value = 0.1447113485548341
$ws.Range("Q6").Value = 33.28329622735656
$ws.Range("R6").Value = 299.549666046209
$ws.Range("S6").Value = 0.04862417841708849
$ws.Range("T6").Value = 0.04862417841708849

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 60.20577233333334
$ws.Range("H7").Value = 180.617317
$ws.Range("I7").Value = 0.3360080526004068
$ws.Range("J7").Value = 0.3360080526004068
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.5834963333333333
$ws.Range("N7").Value = 1.750489
$ws.Range("O7").Value = 0.1527399076504545
$ws.Range("P7").Value = 0.1527399076504546
$ws.Range("Q7").Value = 35.12984740200145
$ws.Range("R7").Value = 316.168626618013
$ws.Range("S7").Value = 0.05132183892399521
$ws.Range("T7").Value = 0.05132183892399522

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 50.58089766666667
$ws.Range("H8").Value = 151.742693
$ws.Range("I8").Value = 0.2822916850839468
$ws.Range("J8").Value = 0.2822916850839468
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.683873666666667
$ws.Range("N8").Value = 8.051621
$ws.Range("O8").Value = 0.7025487437947114
$ws.Range("P8").Value = 0.7025487437947114
$ws.Range("Q8").Value = 135.7527392839281
$ws.Range("R8").Value = 1221.774653555353
$ws.Range("S8").Value = 0.1983236687394191
$ws.Range("T8").Value = 0.1983236687394191

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 50.58089766666667
$ws.Range("H9").Value = 151.742693
$ws.Range("I9").Value = 0.2822916850839468
$ws.Range("J9").Value = 0.2822916850839468
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.5528256666666667
$ws.Range("N9").Value = 1.658477
$ws.Range("O9").Value = 0.1447113485548341
$ws.Range("P9").Value = 0.1447113485548341
$ws.Range("Q9").Value = 27.96241847317344
$ws.Range("R9").Value = 251.661766258561
$ws.Range("S9").Value = 0.04085081043431447
$ws.Range("T9").Value = 0.04085081043431447

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.58089766666667
$ws.Range("H10").Value = 151.742693
$ws.Range("I10").Value = 0.2822916850839468
$ws.Range("J10").Value = 0.2822916850839468
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.5834963333333333
$ws.Range("N10").Value = 1.750489
$ws.Range("O10").Value = 0.1527399076504545
$ws.Range("P10").Value = 0.1527399076504546
$ws.Range("Q10").Value = 29.51376832520855
$ws.Range("R10").Value = 265.623914926877
$ws.Range("S10").Value = 0.04311720591021322
$ws.Range("T10").Value = 0.04311720591021323

